$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC sheet ---
# row 61 (diff hunk @ 3693)
$wsALC.Range("H61").Value = 3466.6667
$wsALC.Range("I61").Value = 200
$wsALC.Range("J61").Value = 10000
$wsALC.Range("K61").Value = 600
$wsALC.Range("L61").Value = 30000
$wsALC.Range("M61").Value = -428
$wsALC.Range("N61").Value = -30344

# row 74 (diff hunk @ 4345)
$wsALC.Range("H74").Value = 2638.4375
$wsALC.Range("I74").Value = 2126.182
$wsALC.Range("J74").Value = 2906.762
$wsALC.Range("K74").Value = 2126.182
$wsALC.Range("L74").Value = 2906.762
$wsALC.Range("M74").Value = -1190.182
$wsALC.Range("N74").Value = -4778.762000000001

# row 77 (diff hunk @ 4498)
$wsALC.Range("H77").Value = 2638.4375
$wsALC.Range("I77").Value = 2126.182
$wsALC.Range("J77").Value = 2906.762
$wsALC.Range("K77").Value = 10630.91
$wsALC.Range("L77").Value = 14533.81
$wsALC.Range("M77").Value = -5950.91
$wsALC.Range("N77").Value = -23893.81

# row 112 (diff hunk @ 6255)
$wsALC.Range("H112").Value = 1858478.9
$wsALC.Range("I112").Value = 934.2857
$wsALC.Range("J112").Value = 2477660.5
$wsALC.Range("K112").Value = 2802.8571
$wsALC.Range("L112").Value = 7432981.5
$wsALC.Range("M112").Value = -1694.8571
$wsALC.Range("N112").Value = -7435197.5

# row 113 (diff hunk @ 6307)
$wsALC.Range("H113").Value = 2276.7856
$wsALC.Range("I113").Value = 1745
$wsALC.Range("J113").Value = 2808.5715
$wsALC.Range("K113").Value = 1745
$wsALC.Range("L113").Value = 2808.5715
$wsALC.Range("M113").Value = 1509
$wsALC.Range("N113").Value = -9316.5715

# row 116 (diff hunk @ 6460)
$wsALC.Range("H116").Value = 3322.682
$wsALC.Range("I116").Value = 2889.2666
$wsALC.Range("K116").Value = 2889.2666
$wsALC.Range("M116").Value = 552.7334000000001

# --- ARM sheet ---
# row 9 (diff hunk @ 8213)
$wsARM.Range("H9").Value = 0
$wsARM.Range("J9").Value = 0
$wsARM.Range("L9").Value = 0
$wsARM.Range("N9").ClearContents()

# row 20 (diff hunk @ 8743)
$wsARM.Range("H20").Value = 0
$wsARM.Range("J20").Value = 0
$wsARM.Range("L20").Value = 0
$wsARM.Range("N20").ClearContents()

# row 52 (diff hunk @ 10317)
$wsARM.Range("H52").Value = 0
$wsARM.Range("J52").Value = 0
$wsARM.Range("L52").Value = 0
$wsARM.Range("N52").ClearContents()

# row 63 (diff hunk @ 10853)
$wsARM.Range("H63").Value = 1000
$wsARM.Range("I63").Value = 1000
$wsARM.Range("J63").Value = 0
$wsARM.Range("K63").Value = 1000
$wsARM.Range("L63").Value = 0
$wsARM.Range("M63").Value = -314
$wsARM.Range("N63").ClearContents()

# row 66 (diff hunk @ 11000)
$wsARM.Range("H66").Value = 1000
$wsARM.Range("I66").Value = 1000
$wsARM.Range("J66").Value = 0
$wsARM.Range("K66").Value = 5000
$wsARM.Range("L66").Value = 0
$wsARM.Range("M66").Value = -1568
$wsARM.Range("N66").ClearContents()

# row 76 (diff hunk @ 11481)
$wsARM.Range("H76").Value = 41519.332
$wsARM.Range("J76").Value = 41519.332
$wsARM.Range("L76").Value = 41519.332
$wsARM.Range("N76").Value = -42195.332

# row 79 (diff hunk @ 11628)
$wsARM.Range("H79").Value = 41519.332
$wsARM.Range("J79").Value = 41519.332
$wsARM.Range("L79").Value = 41519.332
$wsARM.Range("N79").Value = -43859.332

# row 122 (diff hunk @ 13750)
$wsARM.Range("H122").Value = 1862
$wsARM.Range("I122").Value = 0
$wsARM.Range("J122").Value = 1862
$wsARM.Range("K122").Value = 0
$wsARM.Range("L122").Value = 5586
$wsARM.Range("M122").ClearContents()
$wsARM.Range("N122").Value = -10486

# --- BSM sheet ---
# row 105 (diff hunk @ 19895)
$wsBSM.Range("H105").Value = 2303.8125
$wsBSM.Range("I105").Value = 1362.5
$wsBSM.Range("K105").Value = 1362.5
$wsBSM.Range("M105").Value = 384.5

# row 107 (diff hunk @ 19996)
$wsBSM.Range("H107").Value = 1985.8422
$wsBSM.Range("I107").Value = 1607.4828
$wsBSM.Range("J107").Value = 3205
$wsBSM.Range("K107").Value = 1607.4828
$wsBSM.Range("L107").Value = 3205
$wsBSM.Range("M107").Value = 312.5172
$wsBSM.Range("N107").Value = -7045

# --- CRP sheet ---
# row 41 (diff hunk @ 23740)
$wsCRP.Range("H41").Value = 9375
$wsCRP.Range("I41").Value = 3000
$wsCRP.Range("J41").Value = 11500
$wsCRP.Range("K41").Value = 3000
$wsCRP.Range("L41").Value = 11500
$wsCRP.Range("M41").Value = -2572
$wsCRP.Range("N41").Value = -12356

# row 68 (diff hunk @ 25078)
$wsCRP.Range("H68").Value = 20166.666
$wsCRP.Range("I68").Value = 13500
$wsCRP.Range("K68").Value = 13500
$wsCRP.Range("M68").Value = -12751

# row 71 (diff hunk @ 25231)
$wsCRP.Range("H71").Value = 20166.666
$wsCRP.Range("I71").Value = 13500
$wsCRP.Range("K71").Value = 40500
$wsCRP.Range("M71").Value = -36756

# row 74 (diff hunk @ 25384)
$wsCRP.Range("H74").Value = 20157
$wsCRP.Range("J74").Value = 20157
$wsCRP.Range("L74").Value = 20157
$wsCRP.Range("N74").Value = -21905

# row 77 (diff hunk @ 25522)
$wsCRP.Range("H77").Value = 20157
$wsCRP.Range("J77").Value = 20157
$wsCRP.Range("L77").Value = 60471
$wsCRP.Range("N77").Value = -69207

# row 88 (diff hunk @ 26058)
$wsCRP.Range("H88").Value = 16323.2
$wsCRP.Range("J88").Value = 16323.2
$wsCRP.Range("L88").Value = 16323.2
$wsCRP.Range("N88").Value = -17135.2

# row 91 (diff hunk @ 26211)
$wsCRP.Range("H91").Value = 16323.2
$wsCRP.Range("J91").Value = 16323.2
$wsCRP.Range("L91").Value = 16323.2
$wsCRP.Range("N91").Value = -19131.2

# row 99 (diff hunk @ 26606)
$wsCRP.Range("H99").Value = 1570.48
$wsCRP.Range("I99").Value = 1354.375
$wsCRP.Range("J99").Value = 1954.6666
$wsCRP.Range("K99").Value = 1354.375
$wsCRP.Range("L99").Value = 1954.6666
$wsCRP.Range("M99").Value = 143.625
$wsCRP.Range("N99").Value = -4950.6666

# row 107 (diff hunk @ 27004)
$wsCRP.Range("H107").Value = 533.913
$wsCRP.Range("I107").Value = 541.8
$wsCRP.Range("J107").Value = 508.81818
$wsCRP.Range("K107").Value = 541.8
$wsCRP.Range("L107").Value = 508.81818
$wsCRP.Range("M107").Value = 1378.2
$wsCRP.Range("N107").Value = -4348.81818

# row 122 (diff hunk @ 27745)
$wsCRP.Range("H122").Value = 984
$wsCRP.Range("I122").Value = 919.1579
$wsCRP.Range("J122").Value = 1600
$wsCRP.Range("K122").Value = 2757.4737
$wsCRP.Range("L122").Value = 4800
$wsCRP.Range("M122").Value = -307.4737
$wsCRP.Range("N122").Value = -9700

# row 126 (diff hunk @ 27944)
$wsCRP.Range("H126").Value = 1570.48
$wsCRP.Range("I126").Value = 1354.375
$wsCRP.Range("J126").Value = 1954.6666
$wsCRP.Range("K126").Value = 4063.125
$wsCRP.Range("L126").Value = 5863.9998
$wsCRP.Range("M126").Value = -1593.125
$wsCRP.Range("N126").Value = -10803.9998

# --- CUL sheet ---
# row 125 (diff hunk @ 35050)
$wsCUL.Range("H125").Value = 3886.6
$wsCUL.Range("J125").Value = 4433
$wsCUL.Range("L125").Value = 13299
$wsCUL.Range("N125").Value = -23139

# row 130 (diff hunk @ 35304)
$wsCUL.Range("H130").Value = 2203.3333
$wsCUL.Range("J130").Value = 2671.4285
$wsCUL.Range("L130").Value = 8014.2855
$wsCUL.Range("N130").Value = -18054.2855

# --- GSM sheet ---
# row 102 (diff hunk @ 40946)
$wsGSM.Range("H102").Value = 2167.4666
$wsGSM.Range("I102").Value = 2133.7646
$wsGSM.Range("J102").Value = 2211.5386
$wsGSM.Range("K102").Value = 2133.7646
$wsGSM.Range("L102").Value = 2211.5386
$wsGSM.Range("M102").Value = -511.7646
$wsGSM.Range("N102").Value = -5455.5386

# row 132 (diff hunk @ 42416)
$wsGSM.Range("H132").Value = 34936.47
$wsGSM.Range("I132").Value = 54415.906
$wsGSM.Range("J132").Value = 3469.6924
$wsGSM.Range("K132").Value = 163247.718
$wsGSM.Range("L132").Value = 10409.0772
$wsGSM.Range("M132").Value = -160717.718
$wsGSM.Range("N132").Value = -15469.0772

# row 133 (diff hunk @ 42468)
$wsGSM.Range("H133").Value = 33000
$wsGSM.Range("J133").Value = 33000
$wsGSM.Range("L133").Value = 33000
$wsGSM.Range("N133").Value = -43120

# --- LTW sheet ---
# row 7 (diff hunk @ 43254)
$wsLTW.Range("H7").Value = 2036.1333
$wsLTW.Range("I7").Value = 1887.0769
$wsLTW.Range("J7").Value = 3005
$wsLTW.Range("K7").Value = 1887.0769
$wsLTW.Range("L7").Value = 3005
$wsLTW.Range("M7").Value = -1775.0769
$wsLTW.Range("N7").Value = -3229

# row 68 (diff hunk @ 46252)
$wsLTW.Range("H68").Value = 1860
$wsLTW.Range("I68").Value = 1985.7142
$wsLTW.Range("J68").Value = 1750
$wsLTW.Range("K68").Value = 1985.7142
$wsLTW.Range("L68").Value = 1750
$wsLTW.Range("M68").Value = -1236.7142
$wsLTW.Range("N68").Value = -3248

# row 71 (diff hunk @ 46399)
$wsLTW.Range("H71").Value = 1860
$wsLTW.Range("I71").Value = 1985.7142
$wsLTW.Range("J71").Value = 1750
$wsLTW.Range("K71").Value = 9928.571
$wsLTW.Range("L71").Value = 8750
$wsLTW.Range("M71").Value = -6184.571
$wsLTW.Range("N71").Value = -16238

# row 126 (diff hunk @ 49094)
$wsLTW.Range("H126").Value = 2036.1333
$wsLTW.Range("I126").Value = 1887.0769
$wsLTW.Range("J126").Value = 3005
$wsLTW.Range("K126").Value = 5661.2307
$wsLTW.Range("L126").Value = 9015
$wsLTW.Range("M126").Value = -3191.2307
$wsLTW.Range("N126").Value = -13955

# --- WVR sheet ---
# row 107 (diff hunk @ 55132)
$wsWVR.Range("H107").Value = 607.0833
$wsWVR.Range("I107").Value = 527.9
$wsWVR.Range("J107").Value = 1003
$wsWVR.Range("K107").Value = 1583.7
$wsWVR.Range("L107").Value = 3009
$wsWVR.Range("M107").Value = 336.3000000000002
$wsWVR.Range("N107").Value = -6849

# row 122 (diff hunk @ 55864)
$wsWVR.Range("H122").Value = 10001076
$wsWVR.Range("I122").Value = 11765500
$wsWVR.Range("J122").Value = 2670
$wsWVR.Range("K122").Value = 35296500
$wsWVR.Range("L122").Value = 8010
$wsWVR.Range("M122").Value = -35294050
$wsWVR.Range("N122").Value = -12910

